$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with new columns P1 and Q1, copying the header style (s="1") from O1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25 for columns B through Q (adds new columns G, H, P, Q; updates B-F, I-K, N)
$data = @{
    2 = @(24.03912910964759,17.43610907751401,5.645419062495939,6.863109225107139,40.41990741977342,0,2.02613036804859,2.929255597982565,4.251792669778333,27.62635374384049,0,0,8.25274130543689,0,0,0)
    3 = @(22.47261134113976,16.23339414237104,5.574876076283084,6.562816623085172,38.44551530755749,0,2.355117286357314,3.22583832501543,4.251792669778333,26.64381156498273,0,0,7.856814288907563,0,0,0)
    4 = @(21.45801825182291,15.45896180937642,5.528675005221321,6.373439932185669,37.18730761592249,0,2.562988629145515,3.413803134596491,4.251792669778333,26.02100672018133,0,0,7.604309996410807,0,0,0)
    5 = @(21.03140190764024,15.14192871933085,5.504085428723505,6.293521609329384,36.6266948413501,0,2.649863593394134,3.494900808072938,4.251792669778333,25.73369173166343,0,0,7.499425713406922,0,0,0)
    6 = @(20.95985908459775,15.09814327892468,5.493813824529588,6.278419138611069,36.48793808496594,0,2.665162833237697,3.512278863387146,4.251792669778333,25.65045341724782,0,0,7.482264503851145,0,0,0)
    7 = @(21.45252421998708,15.47974802747958,5.511615426951042,6.367640651089719,37.05888121126442,0,2.566264456392319,3.425108981581368,4.251792669778333,25.92193573880895,0,0,7.603959662075227,0,0,0)
    8 = @(23.50999946836863,17.05909813246656,5.59989059134648,6.754569148080214,39.59795798386184,0,2.141082012570588,3.043836104774113,4.251792669778333,27.17077563643697,0,0,8.119532357579198,0,0,0)
    9 = @(27.22110625245902,19.8364747778422,5.779289486034574,7.48434925417068,44.39614458382863,0,1.768110757602306,2.628189607404491,4.251792669778333,29.61770591649753,0,0,9.052315245553363,0,0,0)
    10 = @(29.69407812620388,21.77003456071697,5.813092505016608,7.948533901803495,47.14495020338869,0,2.296829811427553,3.107831974274635,4.251792669778333,30.8795078423273,0,0,9.548846427370863,0,0,0)
    11 = @(30.71714750633855,23.04909056168007,5.215166169136807,7.967684885430201,44.08267507988018,0,3.128811095017249,3.2296556181632,4.251792669778333,27.99836039540484,0,0,8.576659603067302,0,0,0)
    12 = @(31.08140822983099,23.71428086375597,4.835590077963917,8.237095823167987,41.04575072027168,0,4.321438270633167,3.239089003840233,4.251792669778333,25.43010573404221,0,0,7.67729388180924,0,0,0)
    13 = @(30.97248295389151,23.9974909395488,4.601997078627847,8.677262167471365,37.60475593327918,0,5.591488008146638,3.164524792788818,4.251792669778333,22.7850086680612,0,0,6.775448129167177,0,0,0)
    14 = @(30.68160960064628,24.03526655619362,4.532919239349174,9.085263914223495,34.94502909920057,0,6.510603544376936,3.076182043292491,4.251792669778333,20.87856573988034,0,0,6.155560314999028,0,0,0)
    15 = @(30.50797800485661,23.97027341432476,4.53137144537576,9.17761949149271,34.16770191686329,0,6.7178028516996,3.036711014760734,4.251792669778333,20.36926185689898,0,0,6.000882759869356,0,0,0)
    16 = @(29.51677813463343,23.16513016743827,4.560942009403266,8.90380744784593,33.41201322381095,0,6.439084942391636,2.849357922887935,4.251792669778333,20.22511812988133,0,0,5.92842926077746,0,0,0)
    17 = @(28.90093901282754,22.51120415248718,4.60329701756044,8.4392198971971,34.30192828379266,0,5.697235125753479,2.752205321702972,4.251792669778333,21.13429855146545,0,0,6.201873041377087,0,0,0)
    18 = @(28.55352912184775,21.90707190813053,4.744702034893908,7.882809839257665,36.79775442675756,0,4.426981312942706,2.722655203146409,4.251792669778333,23.15402528928243,0,0,6.847632577513944,0,0,0)
    19 = @(28.45402956249922,21.43950799709829,5.050776524327895,7.533587456552659,40.18218848207431,0,3.052870348615724,2.764369529639837,4.251792669778333,25.80708528252396,0,0,7.78614386429331,0,0,0)
    20 = @(29.05999102512012,21.33405619807576,5.75331138391141,7.812185240444131,46.10651996385858,0,2.151606733034478,2.988756946253532,4.251792669778333,30.28036724367185,0,0,9.417088877230281,0,0,0)
    21 = @(30.88141566495933,22.67845571236356,5.902977270865251,8.247045065517062,48.96342631947059,0,2.57897219329035,3.368811507544523,4.251792669778333,31.87792826047457,0,0,10.00128824687754,0,0,0)
    22 = @(32.02286529280498,23.53299585390473,5.979746121177061,8.507958673202372,50.65414323074868,0,2.848636062064426,3.608043894206066,4.251792669778333,32.80195010716557,0,0,10.30046457845435,0,0,0)
    23 = @(31.41905442651385,23.05890991878999,5.957322528078084,8.373583033465874,49.86655221723385,0,2.707131697775672,3.478569147414445,4.251792669778333,32.40368041858036,0,0,10.14064767237077,0,0,0)
    24 = @(29.03499223596726,21.23965733624409,5.838997887442645,7.853723959410925,46.62334538341204,0,2.162083603598564,2.987387392665708,4.251792669778333,30.71184451088924,0,0,9.518946127976799,0,0,0)
    25 = @(26.25472999238668,19.15126690619843,5.704187815245031,7.28271271943919,42.95018278136639,0,1.561645587580376,2.534939558221962,4.251792669778333,28.81829853035774,0,0,8.810049902400866,0,0,0)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    $arr = New-Object "object[,]" 1,16
    for ($i = 0; $i -lt 16; $i++) { $arr[0,$i] = $rowVals[$i] }
    $ws.Range("B$r`:Q$r").Value = $arr
}
